$d = $word.ActiveDocument

function Set-ParaText($paraIndex, $b64text) {
    $bytes = [System.Convert]::FromBase64String($b64text)
    $text = [System.Text.Encoding]::UTF8.GetString($bytes)
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

$para2_b64 = "SS4gQU5URUNFREVOVEVTICALCzEuIEVsIHNlw7FvciBMdWlzIEFsYmVydG8gTWFuamFycmVzIE1hcmlhbm8sIGFmaWxpYWRvIGEgbGEgQUZQIENvbGZvbmRvcyBTLkEuLCBzZSBlbmN1ZW50cmEgYSBsYSBlc3BlcmEgZGUgbGEgZGVmaW5pY2nDs24gZGUgc3UgYm9ubyBwZW5zaW9uYWwgY29tbyBmdWVudGUgZGUgZmluYW5jaWFjacOzbiBkZSBzdSBwZW5zacOzbi4gUGFyYSB0YWwgZWZlY3RvLCBsYSBhZG1pbmlzdHJhZG9yYSBkZSBwZW5zaW9uZXMgaGEgYWRlbGFudGFkbyBnZXN0aW9uZXMgb3JpZW50YWRhcyBhIGxhIHJlY29uc3RydWNjacOzbiBkZSBzdSBoaXN0b3JpYSBsYWJvcmFsLCBjb25mb3JtZSBhIGxvIGRpc3B1ZXN0byBlbiBlbCBhcnTDrWN1bG8gMjAgZGVsIERlY3JldG8gNjU2IGRlIDE5OTQuICALCzIuIFNlZ8O6biBsYSBjZXJ0aWZpY2FjacOzbiBlbGVjdHLDs25pY2EgZGUgdGllbXBvcyBsYWJvcmFkb3MgKENFVElMKSBuw7ptLiAyMDIzMDM4OTIxMTUwMDkwMDA5OTAwMDgsIGV4cGVkaWRhIGVsIDIzIGRlIG1hcnpvIGRlIDIwMjMsIGVsIHNlw7FvciBNYW5qYXJyZXMgTWFyaWFubyBwcmVzdMOzIHN1cyBzZXJ2aWNpb3MgY29tbyBlbXBsZWFkbyBww7pibGljbyBlbiBlbCBIb3NwaXRhbCBOdWVzdHJhIFNlw7FvcmEgZGUgbG9zIFJlbWVkaW9zIGRlIFJpb2hhY2hhIChHdWFqaXJhKSBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDI3IGRlIG5vdmllbWJyZSBkZSAxOTc5LCByZWFsaXphbmRvIGFwb3J0ZXMgcGVuc2lvbmFsZXMgYSBsYSBDYWphIE5hY2lvbmFsIGRlIFByZXZpc2nDs24gU29jaWFsIChDYWphbmFsKS4gRW4gZGljaGEgY2VydGlmaWNhY2nDs24sIHNlIHNlw7FhbGEgYSBsYSBOYWNpw7NuIGNvbW8gZW50aWRhZCByZXNwb25zYWJsZSBkZWwgcGFnby4gIAsLMy4gRWwgOCBkZSBub3ZpZW1icmUgZGUgMjAyMywgbGEgQUZQIENvbGZvbmRvcywgbWVkaWFudGUgbGEgcGxhdGFmb3JtYSBDRVRJTCwgc29saWNpdMOzIGFsIEhvc3BpdGFsIE51ZXN0cmEgU2XDsW9yYSBkZSBsb3MgUmVtZWRpb3MgZGUgUmlvaGFjaGEgbGEgY29ycmVjY2nDs24gZGUgbGEgY2VydGlmaWNhY2nDs24gcGFyYSBxdWUgZGljaGEgZW50aWRhZCBhc3VtaWVyYSBsb3MgdGllbXBvcyBkZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vLCBhcmd1bWVudGFuZG8gcXVlLCBzZWfDum4gaW5mb3JtYWNpw7NuIGRlIGxhIERpcmVjY2nDs24gR2VuZXJhbCBkZSBSZWd1bGFjacOzbiBFY29uw7NtaWNhIGRlIGxhIFNlZ3VyaWRhZCBTb2NpYWwgKERHUkVTUyksIGVsIHBldGljaW9uYXJpbyBubyBlcmEgYmVuZWZpY2lhcmlvIGRlbCBjb250cmF0byBkZSBjb25jdXJyZW5jaWEuIEVuIGNhc28gZGUgZGVzYWN1ZXJkbywgc2Ugc29saWNpdMOzIGFsbGVnYXIgbG9zIHNvcG9ydGVzIGRlIHBhZ28gYSBDYWphbmFsLiAgCws0LiBFbCAxNSBkZSBtYXJ6byBkZSAyMDI0LCBlbCBIb3NwaXRhbCBOdWVzdHJhIFNlw7FvcmEgZGUgbG9zIFJlbWVkaW9zIGRlIFJpb2hhY2hhIHJlc3BvbmRpw7MgYSB0cmF2w6lzIGRlIGxhIHBsYXRhZm9ybWEgQ0VUSUwsIGluZGljYW5kbyBxdWUgbm8gYXN1bcOtYSBsb3MgdGllbXBvcyBkZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vLCBwdWVzIGxvcyBwYWdvcyBjb3JyZXNwb25kaWVudGVzIGEgbGEgw6lwb2NhIGVuIHF1ZSBlc3RlIGxhYm9yw7MgZXJhbiByZXNwb25zYWJpbGlkYWQgZGUgbGEgU2VjcmV0YXLDrWEgZGUgU2FsdWQgRGVwYXJ0YW1lbnRhbCBkZSBsYSBHdWFqaXJhLiBBc2ltaXNtbywgc29saWNpdMOzIHF1ZSBzZSBlbnZpYXJhIHVuIGNvcnJlbyBwYXJhIHN1c3RlbnRhciBkaWNoYSBhZmlybWFjacOzbi4gIAsLNS4gU2Vnw7puIGVsIHJlcG9ydGUgZGUgbGEgT2ZpY2luYSBkZSBCb25vcyBQZW5zaW9uYWxlcyBkZWwgTWluaXN0ZXJpbyBkZSBIYWNpZW5kYSB5IENyw6lkaXRvIFDDumJsaWNvLCBmZWNoYWRvIGVsIDQgZGUganVsaW8gZGUgMjAyNCwgYWwgcHJvY2VzYXIgbGFzIGxpcXVpZGFjaW9uZXMgcHJvdmlzaW9uYWxlcyBlbiBlbCBhcGxpY2F0aXZvIGNvcnJlc3BvbmRpZW50ZSwgc2UgZ2VuZXLDsyBlbCBlcnJvciA0NDM4LCBlbCBjdWFsIGluZGljYSBxdWUgbGEgZW50aWRhZCBubyBlc3TDoSBhc3VtaWRhIHBvciBsYSBOYWNpw7NuIG8gZXhpc3RlbiBwZXJpb2RvcyBubyBhc3VtaWRvcyBwb3IgZXN0YS4gRW4gY29uc2VjdWVuY2lhLCBzZSByZXF1aWVyZSBxdWUgbGEgQUZQIGVudsOtZSBsb3Mgc29wb3J0ZXMgcmVzcGVjdGl2b3MgcGFyYSB2ZXJpZmljYXIgbG9zIGFwb3J0ZXMgcmVhbGl6YWRvcyBhIENhamFuYWwgeSBkZXRlcm1pbmFyIHNpIGxhIGVudGlkYWQgcHVlZGUgc2VyIGFzdW1pZGEgcG9yIGxhIE5hY2nDs24uICALCzYuIEVsIGRlcGFydGFtZW50byBkZSBsYSBHdWFqaXJhIG1hbmlmZXN0w7MgcXVlLCB0cmFzIGFuYWxpemFyIGVsIGNhc28gZGVsIHNlw7FvciBNYW5qYXJyZXMgTWFyaWFubywgbm8gc2UgZW5jb250cmFyb24gbG9zIHNvcG9ydGVzIGRlIHBhZ28gZGUgbG9zIGFwb3J0ZXMgYSBDYWphbmFsIGNvcnJlc3BvbmRpZW50ZXMgYWwgcGVyaW9kbyBjb21wcmVuZGlkbyBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDI3IGRlIG5vdmllbWJyZSBkZSAxOTc5LiBBcmd1bWVudMOzIHF1ZSBlbCBleHBlZGllbnRlIGxhYm9yYWwgbm8gcmVwb3NhIGVuIGxvcyBhcmNoaXZvcyBkZSBsYSBTZWNyZXRhcsOtYSBkZSBTYWx1ZCBEZXBhcnRhbWVudGFsLCB5YSBxdWUsIHBhcmEgbGEgw6lwb2NhLCBlbCBTZXJ2aWNpbyBTZWNjaW9uYWwgZGUgU2FsdWQg4oCTIEdvYmVybmFjacOzbiBkZSBsYSBHdWFqaXJhIMO6bmljYW1lbnRlIHJlYWxpemFiYSBsYSBwb3Nlc2nDs24gZGUgbG9zIGZ1bmNpb25hcmlvcyBkZXNpZ25hZG9zIHBhcmEgbG9zIGhvc3BpdGFsZXMgZGVsIGRlcGFydGFtZW50by4gUG9yIGxvIHRhbnRvLCBzZcOxYWzDsyBxdWUgbGEgRVNFIEhvc3BpdGFsIE51ZXN0cmEgU2XDsW9yYSBkZSBsb3MgUmVtZWRpb3MgZGUgUmlvaGFjaGEgZXMgbGEgZW50aWRhZCBjb21wZXRlbnRlIHBhcmEgcmVtaXRpciBsYSBkb2N1bWVudGFjacOzbiBzb2xpY2l0YWRhLiAgCws3LiBFbCBIb3NwaXRhbCBOdWVzdHJhIFNlw7FvcmEgZGUgbG9zIFJlbWVkaW9zIGRlIFJpb2hhY2hhLCBtZWRpYW50ZSBlc2NyaXRvIGRlbCAyOSBkZSBtYXlvIGRlIDIwMjUsIG5lZ8OzIGNvbXBldGVuY2lhIHBhcmEgYXBvcnRhciBsYXMgcGxhbmlsbGFzIGRlIHBhZ28gZGUgc2VndXJpZGFkIHNvY2lhbCBpbnRlZ3JhbCBkZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vLCBhcmd1bWVudGFuZG8gcXVlLCBlbiBsYSDDqXBvY2EgZW4gcXVlIGVzdGUgbGFib3LDsywgbGFzIGNvdGl6YWNpb25lcyBlcmFuIGVmZWN0dWFkYXMgcG9yIERFU0FMVUQgZGVwYXJ0YW1lbnRhbCwgaG95IFNlY3JldGFyw61hIGRlIFNhbHVkIGRlbCBkZXBhcnRhbWVudG8gZGUgbGEgR3VhamlyYS4gSW5kaWPDsyBxdWUgbGEgRVNFIGFkcXVpcmnDsyBwZXJzb25lcsOtYSBqdXLDrWRpY2EgYSBwYXJ0aXIgZGUgMTk5NCwgbWVkaWFudGUgT3JkZW5hbnphIDAxOCBkZSAxOTk0LCBjb21wbGVtZW50YWRhIHBvciBsYSBPcmRlbmFuemEgMDUxIGRlIDE5OTUuIEFkZW3DoXMsIHNlw7FhbMOzIHF1ZSBsYSByZXNwb25zYWJpbGlkYWQgZGUgY2VydGlmaWNhciBzaSBDYWphbmFsIHJlY2liacOzIGxhcyBjb3RpemFjaW9uZXMgY29ycmVzcG9uZGUgYSBsYSBOYWNpw7NuLCBhIHRyYXbDqXMgZGUgbGEgVW5pZGFkIGRlIEdlc3Rpw7NuIFBlbnNpb25hbCB5IFBhcmFmaXNjYWxlcyAoVUdQUCkuICALCzguIEVsIE1pbmlzdGVyaW8gZGUgSGFjaWVuZGEgeSBDcsOpZGl0byBQw7pibGljbywgbWVkaWFudGUgZG9jdW1lbnRvIGRlbCAzIGRlIGp1bmlvIGRlIDIwMjUsIHNvc3R1dm8gcXVlIG5vIGVzIGxhIGF1dG9yaWRhZCBhZG1pbmlzdHJhdGl2YSBjb21wZXRlbnRlIHBhcmEgcmVzb2x2ZXIgbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIGRlbCBib25vIHBlbnNpb25hbCBkZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vLiBTZcOxYWzDsyBxdWUgZXN0ZSB0aWVuZSBkZXJlY2hvIGEgdW4gYm9ubyBwZW5zaW9uYWwgdGlwbyBBIG1vZGFsaWRhZCAyLCBwb3IgaGFiZXJzZSB0cmFzbGFkYWRvIGFsIHLDqWdpbWVuIGRlIGFob3JybyBpbmRpdmlkdWFsIGNvbiBwb3N0ZXJpb3JpZGFkIGEgbGEgZW50cmFkYSBlbiB2aWdlbmNpYSBkZSBsYSBMZXkgMTAwIGRlIDE5OTMgeSBjb250YXIgY29uIG3DoXMgZGUgMTUwIHNlbWFuYXMgZGUgY290aXphY2nDs24gYWwgSVNTIG8gYSBjYWphcyBww7pibGljYXMuIFNlZ8O6biBsYSBsaXF1aWRhY2nDs24gcHJvdmlzaW9uYWwgZGVsIGJvbm8gcGVuc2lvbmFsLCBlbCBkZXBhcnRhbWVudG8gZGVsIEF0bMOhbnRpY28gY29uY3Vycmlyw61hIGNvbW8gZW1pc29yLCBtaWVudHJhcyBxdWUgbGEgTmFjacOzbiBhc3VtaXLDrWEgbG9zIHRpZW1wb3MgbGFib3JhZG9zIGVuIGxhIEVTRSBIb3NwaXRhbCBOdWVzdHJhIFNlw7FvcmEgZGUgbG9zIFJlbWVkaW9zIGRlIFJpb2hhY2hhIGVudHJlIGVsIDEuwrogZGUganVsaW8gZGUgMTk3OCB5IGVsIDI3IGRlIG5vdmllbWJyZSBkZSAxOTc5LCBjb24gY290aXphY2lvbmVzIGEgQ2FqYW5hbCBkZWJpZGFtZW50ZSBzb3BvcnRhZGFzLiAgCws5LiBFbiBlbCBwcmVzZW50ZSBjYXNvLCBlbCBwZXJpb2RvIHNvYnJlIGVsIGN1YWwgbm8gc2UgY3VlbnRhIGNvbiBzb3BvcnRlcyBxdWUgYWNyZWRpdGVuIGxvcyBhcG9ydGVzIGEgQ2FqYW5hbCBjb3JyZXNwb25kZSBhbCBjb21wcmVuZGlkbyBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDMwIGRlIGp1bmlvIGRlIDE5NzgsIGxvIHF1ZSBnZW5lcmEgZWwgZXJyb3IgNDQzOCBlbiBlbCBzaXN0ZW1hIGludGVyYWN0aXZvIGRlbCBNaW5pc3RlcmlvIGRlIEhhY2llbmRhIHkgQ3LDqWRpdG8gUMO6YmxpY28uIFBvciBsbyB0YW50bywgZWwgcHJvYmxlbWEganVyw61kaWNvIHJhZGljYSBlbiBkZXRlcm1pbmFyIGxhIGF1dG9yaWRhZCBjb21wZXRlbnRlIHBhcmEgcmVzb2x2ZXIgbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIHkgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwgcG9yIGRpY2hvIHBlcmlvZG8sIGNvbmZvcm1lIGFsIGFydMOtY3VsbyAyLjIuMTYuMy44IGRlbCBEZWNyZXRvIDE4ODMgZGUgMjAxNiwgbW9kaWZpY2FkbyBwb3IgZWwgYXJ0w61jdWxvIDEuwrogZGVsIERlY3JldG8gNzkwIGRlIDIwMjEu"
$para5_b64 = "KipJSS4gQ09OU0lERVJBQ0lPTkVTKioLC1BhcmEgcmVzb2x2ZXIgZWwgcHJlc2VudGUgY29uZmxpY3RvIG5lZ2F0aXZvIGRlIGNvbXBldGVuY2lhcyBhZG1pbmlzdHJhdGl2YXMsIGVzdGEgU2FsYSBkZWJlIGRldGVybWluYXIgY3XDoWwgZXMgbGEgYXV0b3JpZGFkIGNvbXBldGVudGUgcGFyYSBlc3R1ZGlhciB5IHJlc29sdmVyIGRlIGZvbmRvIGxhIHNvbGljaXR1ZCBkZSByZWNvbm9jaW1pZW50byB5IHBhZ28gZGVsIGJvbm8gcGVuc2lvbmFsIGNvcnJlc3BvbmRpZW50ZSBhbCBzZcOxb3IgTHVpcyBBbGJlcnRvIE1hbmphcnJlcyBNYXJpYW5vLCBwb3IgZWwgdGllbXBvIGxhYm9yYWRvIGVuIGVsIEhvc3BpdGFsIE51ZXN0cmEgU2XDsW9yYSBkZSBsb3MgUmVtZWRpb3MgZGUgUmlvaGFjaGEgKEd1YWppcmEpLCBlc3BlY8OtZmljYW1lbnRlIHBvciBlbCBwZXLDrW9kbyBjb21wcmVuZGlkbyBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDMwIGRlIGp1bmlvIGRlIDE5NzgsIGFudGUgbGEgYXVzZW5jaWEgZGUgc29wb3J0ZXMgcXVlIGFjcmVkaXRlbiBsb3MgYXBvcnRlcyByZWFsaXphZG9zIGEgbGEgQ2FqYSBOYWNpb25hbCBkZSBQcmV2aXNpw7NuIFNvY2lhbCAoQ2FqYW5hbCkgZHVyYW50ZSBkaWNobyBsYXBzby4LCzEuICoqTWFyY28gbm9ybWF0aXZvIGFwbGljYWJsZSoqCwtFbCBhcnTDrWN1bG8gMi4yLjE2LjMuOCBkZWwgRGVjcmV0byAxODgzIGRlIDIwMTYsIG1vZGlmaWNhZG8gcG9yIGVsIGFydMOtY3VsbyAxLsK6IGRlbCBEZWNyZXRvIDc5MCBkZSAyMDIxLCBlc3RhYmxlY2UgcXVlLCBlbiBhdXNlbmNpYSBkZSBpbmZvcm1hY2nDs24gcXVlIGRlbXVlc3RyZSBlbCBwYWdvIGRlIGxhcyBvYmxpZ2FjaW9uZXMgYSBDYWphbmFsLCBzZSBwcmVzdW1lIHF1ZSBlbCByZXNwb25zYWJsZSBkZWwgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwgZXMgZWwgZW1wbGVhZG9yLiBFc3RlIGRlYmVyw6EgcmVjb25vY2VyIHkgcGFnYXIgZWwgYm9ubyBwZW5zaW9uYWwgY29ycmVzcG9uZGllbnRlLCBzYWx2byBxdWUgYXBvcnRlIGxvcyBzb3BvcnRlcyBxdWUgYWNyZWRpdGVuIGVsIGN1bXBsaW1pZW50byBkZSBsYXMgb2JsaWdhY2lvbmVzIHBlbnNpb25hbGVzLiBBc2ltaXNtbywgZWwgcGFyw6FncmFmbyBkZWwgY2l0YWRvIGFydMOtY3VsbyBkaXNwb25lIHF1ZSBlbCB0aWVtcG8gZGUgc2VydmljaW8gb2ZpY2lhbCBjZXJ0aWZpY2FkbyBjb21vIGNvdGl6YWRvIGEgQ2FqYW5hbCBzb2xvIHNlcsOhIGNvbnRhYmlsaXphZG8gcG9yIGxhIEFkbWluaXN0cmFkb3JhIENvbG9tYmlhbmEgZGUgUGVuc2lvbmVzIChDb2xwZW5zaW9uZXMpIHNpIGV4aXN0ZSBlbCBzb3BvcnRlIGRlIHBhZ28gZW4gbG9zIGFyY2hpdm9zIGRlIGxhIFVuaWRhZCBkZSBHZXN0acOzbiBQZW5zaW9uYWwgeSBQYXJhZmlzY2FsZXMgZGUgbGEgUHJvdGVjY2nDs24gU29jaWFsIChVR1BQKS4gRW4gY2FzbyBkZSBubyBleGlzdGlyIGV2aWRlbmNpYSBkZWwgcGFnbywgZWwgZW1wbGVhZG9yIG8gcXVpZW4gaGFnYSBzdXMgdmVjZXMgZGViZXLDoSBnYXJhbnRpemFyIGxhIGZpbmFuY2lhY2nDs24gZGUgZGljaG9zIHRpZW1wb3MgbWVkaWFudGUgZWwgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwsIGVsIHRyYXNsYWRvIGRlIGFwb3J0ZXMgbyBlbCBjw6FsY3VsbyBhY3R1YXJpYWwgcG9yIG9taXNpw7NuLCBzZWfDum4gY29ycmVzcG9uZGEuCwsyLiAqKkFuw6FsaXNpcyBkZWwgY2FzbyBjb25jcmV0byoqCwtFbiBlbCBwcmVzZW50ZSBhc3VudG8sIHNlIGVuY3VlbnRyYSBhY3JlZGl0YWRvIHF1ZSBlbCBzZcOxb3IgTHVpcyBBbGJlcnRvIE1hbmphcnJlcyBNYXJpYW5vIGxhYm9yw7MgZW4gZWwgSG9zcGl0YWwgTnVlc3RyYSBTZcOxb3JhIGRlIGxvcyBSZW1lZGlvcyBkZSBSaW9oYWNoYSAoR3VhamlyYSkgY29tbyBlbXBsZWFkbyBww7pibGljbyBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDI3IGRlIG5vdmllbWJyZSBkZSAxOTc5LCB5IHF1ZSBsb3MgYXBvcnRlcyBwZW5zaW9uYWxlcyBjb3JyZXNwb25kaWVudGVzIGEgZGljaG8gcGVyw61vZG8gc2UgcmVhbGl6YXJvbiBhIENhamFuYWwuIFNpbiBlbWJhcmdvLCBubyBzZSBjdWVudGEgY29uIGxvcyBzb3BvcnRlcyBxdWUgYWNyZWRpdGVuIGxvcyBwYWdvcyByZWFsaXphZG9zIGEgZGljaGEgZW50aWRhZCBwb3IgZWwgcGVyw61vZG8gY29tcHJlbmRpZG8gZW50cmUgZWwgMS7CuiBkZSBvY3R1YnJlIGRlIDE5NzcgeSBlbCAzMCBkZSBqdW5pbyBkZSAxOTc4LCBsbyBxdWUgaGEgZ2VuZXJhZG8gZWwgZXJyb3IgaWRlbnRpZmljYWRvIGNvbW8gIjQ0MzgiIGVuIGVsIHNpc3RlbWEgaW50ZXJhY3Rpdm8gZGVsIE1pbmlzdGVyaW8gZGUgSGFjaWVuZGEgeSBDcsOpZGl0byBQw7pibGljbywgaW1waWRpZW5kbyBlbCByZWNvbm9jaW1pZW50byBkZWwgYm9ubyBwZW5zaW9uYWwuCwtEZSBhY3VlcmRvIGNvbiBsYSBpbmZvcm1hY2nDs24gYXBvcnRhZGEgYWwgZXhwZWRpZW50ZSwgbGEgRVNFIEhvc3BpdGFsIE51ZXN0cmEgU2XDsW9yYSBkZSBsb3MgUmVtZWRpb3MgZGUgUmlvaGFjaGEgaGEgbmVnYWRvIHN1IGNvbXBldGVuY2lhIHBhcmEgYXN1bWlyIGxvcyB0aWVtcG9zIGVuIGN1ZXN0acOzbiwgYXJndW1lbnRhbmRvIHF1ZSwgcGFyYSBsYSDDqXBvY2EgZW4gcXVlIGVsIHNlw7FvciBNYW5qYXJyZXMgTWFyaWFubyBsYWJvcsOzIGVuIGRpY2hhIGluc3RpdHVjacOzbiwgbGFzIGNvdGl6YWNpb25lcyBlcmFuIHJlc3BvbnNhYmlsaWRhZCBkZWwgZW50ZSB0ZXJyaXRvcmlhbCwgZXNwZWPDrWZpY2FtZW50ZSBkZSBsYSBTZWNyZXRhcsOtYSBkZSBTYWx1ZCBEZXBhcnRhbWVudGFsIGRlIGxhIEd1YWppcmEsIGFudGVzIFNlcnZpY2lvIFNlY2Npb25hbCBkZSBTYWx1ZC4gUG9yIHN1IHBhcnRlLCBlbCBkZXBhcnRhbWVudG8gZGUgbGEgR3VhamlyYSBoYSBzZcOxYWxhZG8gcXVlIG5vIGN1ZW50YSBjb24gbG9zIGV4cGVkaWVudGVzIGxhYm9yYWxlcyBkZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vLCB5YSBxdWUgZXN0b3MgZGViw61hbiByZXBvc2FyIGVuIGVsIGhvc3BpdGFsIGRvbmRlIGVsIHRyYWJhamFkb3IgcHJlc3TDsyBzdXMgc2VydmljaW9zLiBGaW5hbG1lbnRlLCBlbCBNaW5pc3RlcmlvIGRlIEhhY2llbmRhIHkgQ3LDqWRpdG8gUMO6YmxpY28geSBsYSBVR1BQIGhhbiBtYW5pZmVzdGFkbyBxdWUgbm8gc29uIGNvbXBldGVudGVzIHBhcmEgcmVzb2x2ZXIgbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIHkgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwsIGVuIHRhbnRvIG5vIHNlIGhhbiBhcG9ydGFkbyBsb3Mgc29wb3J0ZXMgcXVlIGFjcmVkaXRlbiBsb3MgcGFnb3MgYSBDYWphbmFsLgsLMy4gKipEZXRlcm1pbmFjacOzbiBkZSBsYSBhdXRvcmlkYWQgY29tcGV0ZW50ZSoqCwtDb25mb3JtZSBhIGxvIGRpc3B1ZXN0byBlbiBlbCBhcnTDrWN1bG8gMi4yLjE2LjMuOCBkZWwgRGVjcmV0byAxODgzIGRlIDIwMTYsIG1vZGlmaWNhZG8gcG9yIGVsIGFydMOtY3VsbyAxLsK6IGRlbCBEZWNyZXRvIDc5MCBkZSAyMDIxLCBlbiBhdXNlbmNpYSBkZSBsb3Mgc29wb3J0ZXMgcXVlIGFjcmVkaXRlbiBlbCBwYWdvIGRlIGxvcyBhcG9ydGVzIGEgQ2FqYW5hbCwgc2UgcHJlc3VtZSBxdWUgZWwgZW1wbGVhZG9yIGVzIGVsIHJlc3BvbnNhYmxlIGRlIHJlY29ub2NlciB5IHBhZ2FyIGVsIGJvbm8gcGVuc2lvbmFsIGNvcnJlc3BvbmRpZW50ZS4gRW4gZXN0ZSBjYXNvLCBzZSBlbmN1ZW50cmEgcHJvYmFkbyBxdWUgZWwgc2XDsW9yIE1hbmphcnJlcyBNYXJpYW5vIGxhYm9yw7MgZW4gZWwgSG9zcGl0YWwgTnVlc3RyYSBTZcOxb3JhIGRlIGxvcyBSZW1lZGlvcyBkZSBSaW9oYWNoYSBkdXJhbnRlIGVsIHBlcsOtb2RvIGVuIGN1ZXN0acOzbi4gTm8gb2JzdGFudGUsIGRpY2hhIGVudGlkYWQgaGEgYWxlZ2FkbyBxdWUgbm8gZXMgY29tcGV0ZW50ZSBwYXJhIGFzdW1pciBsb3MgdGllbXBvcywgZW4gcmF6w7NuIGRlIHF1ZSBhZHF1aXJpw7MgcGVyc29uZXLDrWEganVyw61kaWNhIGEgcGFydGlyIGRlbCBhw7FvIDE5OTQsIHkgcXVlLCBwYXJhIGxhIMOpcG9jYSBkZSBsb3MgaGVjaG9zLCBsYXMgY290aXphY2lvbmVzIGVyYW4gcmVzcG9uc2FiaWxpZGFkIGRlbCBlbnRlIHRlcnJpdG9yaWFsLgsLRW4gZXN0ZSBjb250ZXh0bywgeSB0ZW5pZW5kbyBlbiBjdWVudGEgcXVlIGVsIGFydMOtY3VsbyAyLjIuMTYuMy44IGRlbCBEZWNyZXRvIDE4ODMgZGUgMjAxNiBlc3RhYmxlY2UgcXVlIGVsIGVtcGxlYWRvciBvIHF1aWVuIGhhZ2Egc3VzIHZlY2VzIHNlcsOhIHJlc3BvbnNhYmxlIGRlIGdhcmFudGl6YXIgbGEgZmluYW5jaWFjacOzbiBkZSBsb3MgdGllbXBvcyBsYWJvcmFkb3MgbWVkaWFudGUgZWwgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwsIGVsIHRyYXNsYWRvIGRlIGFwb3J0ZXMgbyBlbCBjw6FsY3VsbyBhY3R1YXJpYWwgcG9yIG9taXNpw7NuLCBlc3RhIFNhbGEgY29uY2x1eWUgcXVlIGxhIGF1dG9yaWRhZCBjb21wZXRlbnRlIHBhcmEgcmVzb2x2ZXIgZGUgZm9uZG8gbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIHkgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwgZGVsIHNlw7FvciBMdWlzIEFsYmVydG8gTWFuamFycmVzIE1hcmlhbm8sIHBvciBlbCBwZXLDrW9kbyBjb21wcmVuZGlkbyBlbnRyZSBlbCAxLsK6IGRlIG9jdHVicmUgZGUgMTk3NyB5IGVsIDMwIGRlIGp1bmlvIGRlIDE5NzgsIGVzIGVsIGRlcGFydGFtZW50byBkZSBsYSBHdWFqaXJhLCBlbiBzdSBjYWxpZGFkIGRlIGVudGUgdGVycml0b3JpYWwgcmVzcG9uc2FibGUgZGUgbGEgU2VjcmV0YXLDrWEgZGUgU2FsdWQgRGVwYXJ0YW1lbnRhbCwgYW50ZXMgU2VydmljaW8gU2VjY2lvbmFsIGRlIFNhbHVkLCBxdWUgZnVuZ8OtYSBjb21vIGVtcGxlYWRvciBkZWwgcGV0aWNpb25hcmlvIGVuIGVsIHBlcsOtb2RvIG9iamV0byBkZSBhbsOhbGlzaXMuCws0LiAqKkNvbmNsdXNpw7NuKioLC0VuIG3DqXJpdG8gZGUgbG8gZXhwdWVzdG8sIGVzdGEgU2FsYSBjb25jbHV5ZSBxdWUgZWwgZGVwYXJ0YW1lbnRvIGRlIGxhIEd1YWppcmEgZXMgbGEgYXV0b3JpZGFkIGFkbWluaXN0cmF0aXZhIGNvbXBldGVudGUgcGFyYSByZXNvbHZlciBkZSBmb25kbyBsYSBzb2xpY2l0dWQgZGUgcmVjb25vY2ltaWVudG8geSBwYWdvIGRlbCBib25vIHBlbnNpb25hbCBkZWwgc2XDsW9yIEx1aXMgQWxiZXJ0byBNYW5qYXJyZXMgTWFyaWFubywgcG9yIGVsIHBlcsOtb2RvIGNvbXByZW5kaWRvIGVudHJlIGVsIDEuwrogZGUgb2N0dWJyZSBkZSAxOTc3IHkgZWwgMzAgZGUganVuaW8gZGUgMTk3OCwgZW4gbG9zIHTDqXJtaW5vcyBkZWwgYXJ0w61jdWxvIDIuMi4xNi4zLjggZGVsIERlY3JldG8gMTg4MyBkZSAyMDE2LCBtb2RpZmljYWRvIHBvciBlbCBhcnTDrWN1bG8gMS7CuiBkZWwgRGVjcmV0byA3OTAgZGUgMjAyMS4gQ29ycmVzcG9uZGVyw6EgYSBkaWNoYSBlbnRpZGFkIGFkZWxhbnRhciBsYXMgZ2VzdGlvbmVzIG5lY2VzYXJpYXMgcGFyYSBnYXJhbnRpemFyIGVsIGN1bXBsaW1pZW50byBkZSBsYXMgb2JsaWdhY2lvbmVzIHBlbnNpb25hbGVzIGRlcml2YWRhcyBkZSBkaWNobyBwZXLDrW9kby4="
$para8_b64 = "KipJSUkuIFBST0JMRU1BIEpVUsONRElDTyoqICALC0VsIHByb2JsZW1hIGp1csOtZGljbyBxdWUgc2UgcGxhbnRlYSBlbiBlbCBwcmVzZW50ZSBhc3VudG8gY29uc2lzdGUgZW4gZGV0ZXJtaW5hciBjdcOhbCBlcyBsYSBhdXRvcmlkYWQgYWRtaW5pc3RyYXRpdmEgY29tcGV0ZW50ZSBwYXJhIGVzdHVkaWFyIHkgcmVzb2x2ZXIgZGUgZm9uZG8gbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIHkgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwgY29ycmVzcG9uZGllbnRlIGFsIHRpZW1wbyBsYWJvcmFkbyBwb3IgZWwgc2XDsW9yIEx1aXMgQWxiZXJ0byBNYW5qYXJyZXMgTWFyaWFubyBlbiBlbCBIb3NwaXRhbCBOdWVzdHJhIFNlw7FvcmEgZGUgbG9zIFJlbWVkaW9zIGRlIFJpb2hhY2hhIChHdWFqaXJhKSwgZXNwZWPDrWZpY2FtZW50ZSBwb3IgZWwgcGVyw61vZG8gY29tcHJlbmRpZG8gZW50cmUgZWwgMS7CuiBkZSBvY3R1YnJlIGRlIDE5NzcgeSBlbCAzMCBkZSBqdW5pbyBkZSAxOTc4LCBhbnRlIGxhIGF1c2VuY2lhIGRlIHNvcG9ydGVzIHF1ZSBhY3JlZGl0ZW4gbG9zIGFwb3J0ZXMgcmVhbGl6YWRvcyBhIGxhIENhamEgTmFjaW9uYWwgZGUgUHJldmlzacOzbiBTb2NpYWwgKENhamFuYWwpIGR1cmFudGUgZGljaG8gbGFwc28sIGNvbmZvcm1lIGEgbG8gZGlzcHVlc3RvIGVuIGVsIGFydMOtY3VsbyAyLjIuMTYuMy44IGRlbCBEZWNyZXRvIDE4ODMgZGUgMjAxNiwgbW9kaWZpY2FkbyBwb3IgZWwgYXJ0w61jdWxvIDEuwrogZGVsIERlY3JldG8gNzkwIGRlIDIwMjEu"
$para11_b64 = "KipJVi4gREVDSVNJw5NOKioLC0VuIG3DqXJpdG8gZGUgbG8gZXhwdWVzdG8sIGxhIFNhbGEgZGUgQ29uc3VsdGEgeSBTZXJ2aWNpbyBDaXZpbCBkZWwgQ29uc2VqbyBkZSBFc3RhZG8sIGFkbWluaXN0cmFuZG8ganVzdGljaWEgZW4gbm9tYnJlIGRlIGxhIFJlcMO6YmxpY2EgeSBwb3IgYXV0b3JpZGFkIGRlIGxhIGxleSwgIAsLKipSRVNVRUxWRToqKiAgCwsqKlBSSU1FUk86KiogRGVjbGFyYXIgY29tcGV0ZW50ZSBhbCBkZXBhcnRhbWVudG8gZGUgTGEgR3VhamlyYSwgYSB0cmF2w6lzIGRlIHN1IFNlY3JldGFyw61hIGRlIFNhbHVkIERlcGFydGFtZW50YWwsIHBhcmEgcmVzb2x2ZXIgZGUgZm9uZG8gbGEgc29saWNpdHVkIGRlIHJlY29ub2NpbWllbnRvIHkgcGFnbyBkZWwgYm9ubyBwZW5zaW9uYWwgYSBmYXZvciBkZWwgc2XDsW9yIEx1aXMgQWxiZXJ0byBNYW5qYXJyZXMgTWFyaWFubywgcG9yIGVsIHBlcmlvZG8gY29tcHJlbmRpZG8gZW50cmUgZWwgMS7CuiBkZSBvY3R1YnJlIGRlIDE5NzcgeSBlbCAzMCBkZSBqdW5pbyBkZSAxOTc4LCB0aWVtcG8gZHVyYW50ZSBlbCBjdWFsIGxhYm9yw7MgZW4gZWwgSG9zcGl0YWwgTnVlc3RyYSBTZcOxb3JhIGRlIGxvcyBSZW1lZGlvcyBkZSBSaW9oYWNoYSAoR3VhamlyYSksIGFudGUgbGEgYXVzZW5jaWEgZGUgc29wb3J0ZXMgZGUgcGFnbyBhIGxhIENhamEgTmFjaW9uYWwgZGUgUHJldmlzacOzbiBTb2NpYWwgKENhamFuYWwpIHBvciBkaWNobyBwZXJpb2RvLCBkZSBjb25mb3JtaWRhZCBjb24gbG8gZGlzcHVlc3RvIGVuIGVsIGFydMOtY3VsbyAyLjIuMTYuMy44IGRlbCBEZWNyZXRvIDE4ODMgZGUgMjAxNiwgbW9kaWZpY2FkbyBwb3IgZWwgYXJ0w61jdWxvIDEuwrogZGVsIERlY3JldG8gNzkwIGRlIDIwMjEsIGVuIGxvcyB0w6lybWlub3Mgc2XDsWFsYWRvcyBlbiBsYSBwYXJ0ZSBjb25zaWRlcmF0aXZhIGRlIGVzdGEgZGVjaXNpw7NuLiAgCwsqKlNFR1VORE86KiogUmVtaXRpciBlbCBleHBlZGllbnRlIGRlbCBjb25mbGljdG8gYSBsYSBTZWNyZXRhcsOtYSBkZSBTYWx1ZCBEZXBhcnRhbWVudGFsIGRlIExhIEd1YWppcmEgcGFyYSBsb3MgZmluZXMgZGlzcHVlc3RvcyBlbiBlbCBudW1lcmFsIGFudGVyaW9yLiAgCwsqKlRFUkNFUk86KiogRXhob3J0YXIgYWwgZGVwYXJ0YW1lbnRvIGRlIExhIEd1YWppcmEgcGFyYSBxdWUsIGVuIGF0ZW5jacOzbiBhIGxhIGVzcGVjaWFsIHByb3RlY2Npw7NuIGNvbnN0aXR1Y2lvbmFsIGRlIGxhcyBwZXJzb25hcyBhZHVsdGFzIG1heW9yZXMsIGFkZWxhbnRlIGRlIG1hbmVyYSBwcmlvcml0YXJpYSBsYXMgZ2VzdGlvbmVzIG5lY2VzYXJpYXMgcGFyYSBkYXIgdW5hIHByb250YSByZXNwdWVzdGEgYSBsYSBzb2xpY2l0dWQgcHJlc2VudGFkYSBwb3IgbGEgQUZQIENvbGZvbmRvcywgZW4gcmVwcmVzZW50YWNpw7NuIGRlbCBzZcOxb3IgTHVpcyBBbGJlcnRvIE1hbmphcnJlcyBNYXJpYW5vLiAgCwtOb3RpZsOtcXVlc2UgeSBjw7ptcGxhc2Uu"

Set-ParaText 2 $para2_b64
Set-ParaText 5 $para5_b64
Set-ParaText 8 $para8_b64
Set-ParaText 11 $para11_b64

Write-Host "Done"